# Bill of Materials update: the connector's JLCPCB part number changed
# (new SMA edge-launch connector part), so the "JLCPCB Part #" cell for
# designator J2 (row 2, column D) needs its text updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("D2")
$cell.Value = "C1509219"

# Re-apply the cell's number format so the style record is refreshed
# along with the new part number (this is what produced the extra,
# near-identical style entry seen in the saved workbook).
$cell.NumberFormat = "General"
